{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the documented change:\n//  1. Removes proofErr (spell/grammar check) wrapper markers around\n//     \"get_command\" in the \"main loop\" paragraph (text unchanged).\n//  2. Merges the \"The function named / get_command / contains a while\n//     loop...\" paragraph into a single lead-in run (dropping proofErr),\n//     and collapses the trailing \"(\" + \").\" runs (with the gramStart/\n//     gramEnd proofErr between them) into a single \"().\" run.\n//  3. Splits the \"Each character ... cmdbuffer ...\" paragraph: keeps its\n//     text (minus the proofErr around \"cmdbuffer\") in place, and adds a\n//     brand-new paragraph after it with the new \"array of strings\"\n//     commentary.\n//\n// Implementation note: plain insertText()/Range.text edits in this\n// engine do not clear pre-existing <w:proofErr/> sentinels that wrap a\n// run (they survive even when the run's own text is changed), and a\n// whole-paragraph \"Replace\" via insertText() also collapses the\n// paragraph down to a single merged run (losing the multi-run shape the\n// canonical OOXML expects). To reproduce the exact target markup we\n// instead rebuild each affected paragraph's Range using\n// Range.insertOoxml(..., \"Replace\") with hand-written <w:p>/<w:r> markup\n// -- this lets us control run boundaries precisely and omits the\n// <w:proofErr/> elements outright (equivalent to Word's own behavior of\n// dropping stale proofing sentinels once the reviewed text is edited).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst OOXML_HEADER = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>`;\nconst OOXML_FOOTER = `</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nfunction wrapOoxml(bodyXml) {\n  return OOXML_HEADER + bodyXml + OOXML_FOOTER;\n}\n\n// Locate the three paragraphs we need to touch by matching their\n// (pre-edit) text content, rather than hard-coded indices, so the\n// script is resilient to unrelated paragraphs shifting around.\nparagraphs.items.forEach((p) => {});\n\nlet mainLoopPara = null;\nlet functionNamedPara = null;\nlet eachCharPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"The main loop uses the function named\") === 0) {\n    mainLoopPara = p;\n  } else if (t.indexOf(\"The function named\") === 0 && t.indexOf(\"contains a while loop\") !== -1) {\n    functionNamedPara = p;\n  } else if (t.indexOf(\"Each character from the serial port\") === 0) {\n    eachCharPara = p;\n  }\n}\n\nif (!mainLoopPara || !functionNamedPara || !eachCharPara) {\n  throw new Error(\"Could not locate expected paragraphs to edit.\");\n}\n\n// 1) \"The main loop uses the function named get_command to retrieve...\"\n//    -- same three runs, proofErr markers removed.\n{\n  const xml = wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">The main loop uses the function named </w:t></w:r>' +\n      \"<w:r><w:t>get_command</w:t></w:r>\" +\n      '<w:r><w:t xml:space=\"preserve\"> to retrieve characters from the serial port.</w:t></w:r>' +\n      \"</w:p>\"\n  );\n  mainLoopPara.getRange().insertOoxml(xml, \"Replace\");\n  await context.sync();\n}\n\n// 2) \"The function named get_command contains a while loop...\" -- lead-in\n//    text merged into one run, \"MSerial.available\" kept as its own run,\n//    and the trailing \"(\" + \").\" collapsed into a single \"().\" run.\n{\n  const xml = wrapOoxml(\n    \"<w:p>\" +\n      '<w:r><w:t xml:space=\"preserve\">The function named get_command contains a while loop in which the serial is queried via a call to </w:t></w:r>' +\n      \"<w:r><w:t>MSerial.available</w:t></w:r>\" +\n      \"<w:r><w:t>().</w:t></w:r>\" +\n      \"</w:p>\"\n  );\n  functionNamedPara.getRange().insertOoxml(xml, \"Replace\");\n  await context.sync();\n}\n\n// 3) \"Each character... cmdbuffer...\" -- proofErr around \"cmdbuffer\"\n//    removed (text/run shape otherwise unchanged), followed by a new\n//    paragraph with the added commentary.\n{\n  const xml = wrapOoxml(\n    \"<w:p>\" +\n      \"<w:r><w:t>Each character from the serial port is saved to the buffer named cmdbuffer. This buffer is a 2 dimensional array: 4(</w:t></w:r>\" +\n      \"<w:r><w:t>BUFSIZE</w:t></w:r>\" +\n      \"<w:r><w:t>) x 96(</w:t></w:r>\" +\n      \"<w:r><w:t>MAX_CMD_SIZE</w:t></w:r>\" +\n      \"<w:r><w:t>).</w:t></w:r>\" +\n      \"</w:p>\" +\n      \"<w:p>\" +\n      \"<w:r><w:t>The cmdbuffer is an array of strings! There are a maximum of 4 strings that can be processed at a time. The Marlin firmware will retrieve 4 strings from the serial stream.</w:t></w:r>\" +\n      \"</w:p>\"\n  );\n  eachCharPara.getRange().insertOoxml(xml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Applies the documented change:\n#  1. Removes proofErr (spell/grammar check) wrapper markers around\n#     \"get_command\" in the \"main loop\" paragraph (text unchanged).\n#  2. Merges the \"The function named / get_command / contains a while\n#     loop...\" paragraph into a single lead-in run (dropping proofErr),\n#     and collapses the trailing \"(\" + \").\" runs (with the gramStart/\n#     gramEnd proofErr between them) into a single \"().\" run.\n#  3. Splits the \"Each character ... cmdbuffer ...\" paragraph: keeps its\n#     text (minus the proofErr around \"cmdbuffer\") in place, and adds a\n#     brand-new paragraph after it with the new \"array of strings\"\n#     commentary.\n#\n# Implementation note: plain Range.Text assignment / Find-Replace edits\n# in this engine do not clear pre-existing <w:proofErr/> sentinels that\n# wrap a run (they survive even when the run's own text is changed), and\n# naive whole-paragraph text replacement also collapses the paragraph\n# down to a single merged run (losing the multi-run shape the canonical\n# OOXML expects). To reproduce the exact target markup we instead\n# rebuild each affected paragraph's Range using Range.InsertXML(...)\n# with hand-written <w:p>/<w:r> markup -- this lets us control run\n# boundaries precisely and omits the <w:proofErr/> elements outright\n# (equivalent to Word's own behavior of dropping stale proofing\n# sentinels once the reviewed text is edited).\n\n$d = $word.ActiveDocument\n\n$xmlHeader = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Locate the three paragraphs we need to touch by matching their\n# (pre-edit) text content, rather than a hard-coded index, so the\n# script is resilient to unrelated paragraphs shifting around.\n$mainLoopPara = $null\n$functionNamedPara = $null\n$eachCharPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.StartsWith(\"The main loop uses the function named\")) {\n        $mainLoopPara = $p\n    } elseif ($t.StartsWith(\"The function named\") -and $t.Contains(\"contains a while loop\")) {\n        $functionNamedPara = $p\n    } elseif ($t.StartsWith(\"Each character from the serial port\")) {\n        $eachCharPara = $p\n    }\n}\n\nif ($mainLoopPara -eq $null -or $functionNamedPara -eq $null -or $eachCharPara -eq $null) {\n    throw \"Could not locate expected paragraphs to edit.\"\n}\n\n# 1) \"The main loop uses the function named get_command to retrieve...\"\n#    -- same three runs, proofErr markers removed.\n$xml1 = $xmlHeader + '<w:p><w:r><w:t xml:space=\"preserve\">The main loop uses the function named </w:t></w:r><w:r><w:t>get_command</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to retrieve characters from the serial port.</w:t></w:r></w:p>' + $xmlFooter\n$mainLoopPara.Range.InsertXML($xml1)\n\n# 2) \"The function named get_command contains a while loop...\" -- lead-in\n#    text merged into one run, \"MSerial.available\" kept as its own run,\n#    and the trailing \"(\" + \").\" collapsed into a single \"().\" run.\n$xml2 = $xmlHeader + '<w:p><w:r><w:t xml:space=\"preserve\">The function named get_command contains a while loop in which the serial is queried via a call to </w:t></w:r><w:r><w:t>MSerial.available</w:t></w:r><w:r><w:t>().</w:t></w:r></w:p>' + $xmlFooter\n$functionNamedPara.Range.InsertXML($xml2)\n\n# 3) \"Each character... cmdbuffer...\" -- proofErr around \"cmdbuffer\"\n#    removed (text/run shape otherwise unchanged), followed by a new\n#    paragraph with the added commentary.\n$xml3 = $xmlHeader + '<w:p><w:r><w:t>Each character from the serial port is saved to the buffer named cmdbuffer. This buffer is a 2 dimensional array: 4(</w:t></w:r><w:r><w:t>BUFSIZE</w:t></w:r><w:r><w:t>) x 96(</w:t></w:r><w:r><w:t>MAX_CMD_SIZE</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p><w:p><w:r><w:t>The cmdbuffer is an array of strings! There are a maximum of 4 strings that can be processed at a time. The Marlin firmware will retrieve 4 strings from the serial stream.</w:t></w:r></w:p>' + $xmlFooter\n$eachCharPara.Range.InsertXML($xml3)\n\nWrite-Output \"done\"\n"}
